$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update betting odds values for rows 8, 9, 10, 11, 12, 15, 16
# as per the Jogos_da_Semana_FlashScore_2024-09-24 data refresh

$ws.Range("M8").Value = 1.04
$ws.Range("N8").Value = 13
$ws.Range("G9").Value = 1.9
$ws.Range("I9").Value = 3.8
$ws.Range("J9").Value = 2.6
$ws.Range("L9").Value = 4.5
$ws.Range("S9").Value = 1.4
$ws.Range("T9").Value = 2.75
$ws.Range("U9").Value = 1.83
$ws.Range("V9").Value = 1.83
$ws.Range("W9").Value = 7
$ws.Range("X9").Value = 8.5
$ws.Range("Y9").Value = 9
$ws.Range("AA9").Value = 17
$ws.Range("AG9").Value = 301
$ws.Range("AI9").Value = 19
$ws.Range("AJ9").Value = 13
$ws.Range("AN9").Value = 4
$ws.Range("AT9").Value = 2.75
$ws.Range("AY9").Value = 21
$ws.Range("AZ9").Value = 29
$ws.Range("G10").Value = 2.3
$ws.Range("H10").Value = 2.8
$ws.Range("I10").Value = 3.7
$ws.Range("J10").Value = 3.25
$ws.Range("K10").Value = 1.8
$ws.Range("L10").Value = 4.5
$ws.Range("M10").Value = 1.14
$ws.Range("N10").Value = 5.5
$ws.Range("U10").Value = 2.5
$ws.Range("V10").Value = 1.5
$ws.Range("W10").Value = 5
$ws.Range("X10").Value = 9
$ws.Range("Z10").Value = 21
$ws.Range("AA10").Value = 26
$ws.Range("AH10").Value = 7
$ws.Range("AI10").Value = 15
$ws.Range("AJ10").Value = 15
$ws.Range("AK10").Value = 41
$ws.Range("AL10").Value = 41
$ws.Range("AM10").Value = 51
$ws.Range("AN10").Value = 4
$ws.Range("AO10").Value = 15
$ws.Range("AQ10").Value = 51
$ws.Range("AU10").Value = 10
$ws.Range("AX10").Value = 5.5
$ws.Range("AY10").Value = 23
$ws.Range("BA10").Value = 81
$ws.Range("BB10").Value = 151
$ws.Range("I11").Value = 3.7
$ws.Range("K11").Value = 1.83
$ws.Range("L11").Value = 4.75
$ws.Range("M11").Value = 1.13
$ws.Range("N11").Value = 6
$ws.Range("O11").Value = 1.62
$ws.Range("P11").Value = 2.2
$ws.Range("Q11").Value = 2.88
$ws.Range("R11").Value = 1.4
$ws.Range("S11").Value = 1.67
$ws.Range("T11").Value = 2.1
$ws.Range("U11").Value = 2.5
$ws.Range("V11").Value = 1.5
$ws.Range("W11").Value = 5
$ws.Range("X11").Value = 8.5
$ws.Range("Y11").Value = 11
$ws.Range("AC11").Value = 5.5
$ws.Range("AE11").Value = 23
$ws.Range("AF11").Value = 101
$ws.Range("AH11").Value = 7
$ws.Range("AI11").Value = 17
$ws.Range("AJ11").Value = 15
$ws.Range("AP11").Value = 34
$ws.Range("AR11").Value = 101
$ws.Range("AS11").Value = 351
$ws.Range("AT11").Value = 2.1
$ws.Range("AU11").Value = 10
$ws.Range("BB11").Value = 151
$ws.Range("G12").Value = 3.9
$ws.Range("K12").Value = 1.95
$ws.Range("Q12").Value = 2.4
$ws.Range("R12").Value = 1.53
$ws.Range("S12").Value = 1.53
$ws.Range("T12").Value = 2.38
$ws.Range("U12").Value = 2.1
$ws.Range("V12").Value = 1.67
$ws.Range("X12").Value = 17
$ws.Range("AL12").Value = 21
$ws.Range("AR12").Value = 126
$ws.Range("AT12").Value = 2.38
$ws.Range("BC12").Value = 251
$ws.Range("G15").Value = 2.6
$ws.Range("H15").Value = 3.1
$ws.Range("I15").Value = 2.75
$ws.Range("J15").Value = 3.25
$ws.Range("L15").Value = 3.5
$ws.Range("Q15").Value = 2.2
$ws.Range("R15").Value = 1.65
$ws.Range("S15").Value = 1.5
$ws.Range("T15").Value = 2.5
$ws.Range("X15").Value = 12
$ws.Range("Z15").Value = 26
$ws.Range("AA15").Value = 23
$ws.Range("AC15").Value = 8
$ws.Range("AG15").Value = 351
$ws.Range("AH15").Value = 8
$ws.Range("AI15").Value = 13
$ws.Range("AL15").Value = 23
$ws.Range("AN15").Value = 4.5
$ws.Range("AO15").Value = 15
$ws.Range("AQ15").Value = 51
$ws.Range("AR15").Value = 81
$ws.Range("AT15").Value = 2.5
$ws.Range("AU15").Value = 8.5
$ws.Range("AZ15").Value = 26
$ws.Range("G16").Value = 1.3
$ws.Range("H16").Value = 5
$ws.Range("I16").Value = 10
$ws.Range("J16").Value = 1.8
$ws.Range("L16").Value = 8
$ws.Range("M16").Value = 1.03
$ws.Range("N16").Value = 15
$ws.Range("Q16").Value = 1.65
$ws.Range("R16").Value = 2.2
$ws.Range("W16").Value = 7
$ws.Range("Y16").Value = 8.5
$ws.Range("Z16").Value = 8
$ws.Range("AB16").Value = 29
$ws.Range("AC16").Value = 12
$ws.Range("AH16").Value = 26
$ws.Range("AI16").Value = 51
$ws.Range("AJ16").Value = 29
$ws.Range("AK16").Value = 126
$ws.Range("AL16").Value = 67
$ws.Range("AM16").Value = 67
$ws.Range("AN16").Value = 3.25
$ws.Range("AO16").Value = 6
$ws.Range("AP16").Value = 19
$ws.Range("AU16").Value = 9.5
$ws.Range("AV16").Value = 67
$ws.Range("AX16").Value = 9.5
$ws.Range("BA16").Value = 201
$ws.Range("BB16").Value = 201
$ws.Range("BC16").Value = 351
